$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3496450.76
$ws.Range("C9").Value = 546385.24
$ws.Range("D9").Value = 4042836
$ws.Range("E9").Value = 13.51489993657917
$ws.Range("F9").Value = 86.48510006342082
$ws.Range("G9").Value = -47.19459095640611
$ws.Range("H9").Value = -36.85893220715074
$ws.Range("I9").Value = 35092
$ws.Range("J9").Value = 1490
$ws.Range("K9").Value = 36582
$ws.Range("L9").Value = 25250
$ws.Range("M9").Value = 160.1123168316832
$ws.Range("N9").Value = 9.311939825334491
